$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (replaces the previous Player/Position/Team listing in A2:C18)
$players = @(
    "Damian Lillard",
    "Bilal Coulibaly",
    "Kentavious Caldwell-Pope",
    "Malik Monk",
    "Cameron Johnson",
    "Julius Randle",
    "Brandon Miller",
    "Herbert Jones",
    "Bam Adebayo",
    "Anthony Davis",
    "Nick Richards",
    "Kelly Olynyk",
    "Isaiah Hartenstein",
    "Cade Cunningham",
    "Derrick White",
    "Brandon Ingram",
    "LaMelo Ball"
)

$positions = @(
    "PG",
    "SG,SF",
    "SG,SF",
    "SG,SF",
    "SF,PF",
    "PF",
    "SG,SF",
    "SF,PF",
    "C",
    "PF,C",
    "C",
    "C",
    "C",
    "PG,SG",
    "PG,SG",
    "SG,SF,PF",
    "PG,SG"
)

$teams = @(
    "Milwaukee Bucks",
    "Washington Wizards",
    "Orlando Magic",
    "Sacramento Kings",
    "Brooklyn Nets",
    "Minnesota Timberwolves",
    "Charlotte Hornets",
    "New Orleans Pelicans",
    "Miami Heat",
    "Los Angeles Lakers",
    "Charlotte Hornets",
    "Toronto Raptors",
    "Oklahoma City Thunder",
    "Detroit Pistons",
    "Boston Celtics",
    "New Orleans Pelicans",
    "Charlotte Hornets"
)

for ($i = 0; $i -lt $players.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $players[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
